# Hands-on Intro to Node.js - minor update
#  1. Bump the cached "datetimeFigureOut" footer field from 1/12/15 to
#     1/13/15 on the slide master and every slide layout (mirrors what
#     PowerPoint does when it re-caches the auto date field on save).
#  2. Remove the ad-hoc "TextBox 3" (Wi-fi / login / password) shape
#     that was added to slide 1.

$p = $ppt.ActivePresentation

$oldDate = "1/12/15"
$newDate = "1/13/15"
$ppPlaceholderDate = 16

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ((-not $isDatePlaceholder) -and ($shp.Name -like "Date Placeholder*")) {
            $isDatePlaceholder = $true
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a. Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# 1b. Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 2. Drop the Wi-fi / login / password textbox added to slide 1.
$slide1 = $p.Slides.Item(1)
for ($i = $slide1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 3") {
        $shp.Delete()
    }
}
